$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A69").Value = "Golang + MQTT Developer"
$ws.Range("B69").Value = "https://www.dice.com/job-detail/e0364328-651a-4abb-be87-dbaaec2ea67c"
$ws.Range("C69").Value = "Atlanta, Georgia"
$ws.Range("D69").Value = "Third Party, Contract"
$ws.Range("E69").Value = "`$73 - `$83 per hour"
$ws.Range("F69").Value = "STAND 8"

$ws.Range("A70").Value = "Backend Software Engineer"
$ws.Range("B70").Value = "https://www.dice.com/job-detail/028c9435-f21f-456c-a58b-3f891008f363"
$ws.Range("C70").Value = "Hybrid in Redmond, Washington"
$ws.Range("D70").Value = "Third Party, Contract"
$ws.Range("E70").Value = "Depends on Experience"
$ws.Range("F70").Value = "Technovision, Inc."

$ws.Range("A71").Value = "AWS or Google Cloud Platform Administrator, an architect, or an Admin.|| Onsite - Hybrid at Mclean VA / Santa Monica, CA || Must have Linkedin and 14+ years of exp.||"
$ws.Range("B71").Value = "https://www.dice.com/job-detail/6dbf6404-619e-469d-820b-682c871046f4"
$ws.Range("C71").Value = "Hybrid in McLean, Virginia"
$ws.Range("D71").Value = "Contract, Third Party"
$ws.Range("E71").Value = "Up to `$60"
$ws.Range("F71").Value = "Tri-Force Consulting Services Inc"
